$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy style from existing header cell (H1) to new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data values for columns I (I0) and J (IF), rows 2-28
$data = @(
    @(7, 8),
    @(5, 7),
    @(5, 7),
    @(9, 9),
    @(7, 8),
    @(7, 9),
    @(9, 9),
    @(6, 7),
    @(4, 8),
    @(12, 12),
    @(7, 7),
    @(3, 6),
    @(6, 8),
    @(9, 9),
    @(5, 7),
    @(7, 8),
    @(7, 8),
    @(6, 8),
    @(6, 8),
    @(3, 8),
    @(4, 6),
    @(7, 7),
    @(7, 7),
    @(1, 3),
    @(1, 3),
    @(5, 7),
    @(1, 2)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
